# Applies the outputs-HGR-r202-archive3/g__UMGS75.xlsx edit:
#  - header row (row 1) columns C/D/E are cycled: C<-"prediction", D<-"rejection-f", E<-"max"
#  - data rows 2-4: column C becomes the species string (same text as column D),
#    column D is unchanged, and column E becomes the numeric value 1.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header row
$ws.Range("C1").Value = "prediction"
$ws.Range("D1").Value = "rejection-f"
$ws.Range("E1").Value = "max"

# Data rows
$ws.Range("C2").Value = "s__UMGS75 sp900538885"
$ws.Range("D2").Value = "s__UMGS75 sp900538885"
$ws.Range("E2").Value = 1

$ws.Range("C3").Value = "s__UMGS75 sp900538885"
$ws.Range("D3").Value = "s__UMGS75 sp900538885"
$ws.Range("E3").Value = 1

$ws.Range("C4").Value = "s__UMGS75 sp900538885"
$ws.Range("D4").Value = "s__UMGS75 sp900538885"
$ws.Range("E4").Value = 1
